$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list price (D) and 1h volume change (E) columns
$ws.Range("D2").Value = "'61.952.11"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("D3").Value = "'3.420.13"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'406.29"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").Value = "'132.51"
$ws.Range("E6").Value = "  +3.56%  "
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.690"
$ws.Range("E9").Value = "  +3.23%  "
$ws.Range("D10").Value = "'0.132"
$ws.Range("E10").Value = "  +4.74%  "
$ws.Range("D11").Value = "'41.99"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("D13").Value = "'19.94"
$ws.Range("E13").Value = "  +1.85%  "
$ws.Range("D14").Value = "'8.43"
$ws.Range("E14").Value = "  -1.04%  "
$ws.Range("D15").Value = "'3.412.05"
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").Value = "'11.72"
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("D17").Value = "'61.929.85"
$ws.Range("E17").Value = "  +1.66%  "
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("D19").Value = "'0.0000146"
$ws.Range("E19").Value = "  +10.14%  "
$ws.Range("D20").Value = "'3.17"
$ws.Range("E20").Value = "  -1.64%  "
$ws.Range("D21").Value = "'83.89"
$ws.Range("E21").Value = "  +2.12%  "
$ws.Range("D22").Value = "'313.65"
$ws.Range("E22").Value = "  +2.53%  "
$ws.Range("E23").Value = "  -0.72%  "
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("D25").Value = "'4.75"
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("D26").Value = "'29.67"
$ws.Range("E26").Value = "  +0.67%  "
$ws.Range("E27").Value = "  +6.84%  "
$ws.Range("D28").Value = "'8.14"
$ws.Range("E28").Value = "  -5.54%  "
$ws.Range("D29").Value = "'2.75"
$ws.Range("E29").Value = "  +7.07%  "
$ws.Range("D30").Value = "'0.174"
$ws.Range("E30").Value = "  +0.81%  "
$ws.Range("D31").Value = "'44.01"
$ws.Range("E31").Value = "  +2.95%  "
$ws.Range("E32").Value = "  +0.69%  "
$ws.Range("D33").Value = "'11.33"
$ws.Range("E33").Value = "  -2.60%  "
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("E36").Value = "  -0.80%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").Value = "'3.01"
$ws.Range("E38").Value = "  +1.60%  "
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("D40").Value = "'0.316"
$ws.Range("E40").Value = "  +11.92%  "
$ws.Range("D41").Value = "'140.26"
$ws.Range("E41").Value = "  +3.93%  "
$ws.Range("E43").Value = "  +1.17%  "
$ws.Range("D44").Value = "'3.95"
$ws.Range("E44").Value = "  +1.29%  "
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("D47").Value = "'21.39"
$ws.Range("E47").Value = "  -0.76%  "
$ws.Range("D48").Value = "'2.105.28"
$ws.Range("E48").Value = "  -1.43%  "
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("E50").Value = "  +1.69%  "
$ws.Range("D51").Value = "'1.72"
$ws.Range("E51").Value = "  +18.50%  "
